$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 158, shifting existing rows 158:288 down to 160:290
$ws.Rows("158:159").Insert()

# Populate the new row 158 (Primera) with fresh weekly data
$ws.Range("A158").Value = 8
$ws.Range("B158").Value = "Terminal La Palmera de La Serena"
$ws.Range("C158").Value = "Coquimbo"
$ws.Range("D158").Value = 44554
$ws.Range("E158").Value = 4
$ws.Range("F158").Value = 100112017
$ws.Range("G158").Value = "Apio"
$ws.Range("H158").Value = "Americana (o)"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 2400
$ws.Range("K158").Value = 8000
$ws.Range("L158").Value = 9000
$ws.Range("M158").Value = 8500
$ws.Range("N158").Value = "`$/docena de matas"
$ws.Range("O158").Value = "Provincia del Elquí"
$ws.Range("P158").Value = 1417
$ws.Range("Q158").Value = 6
$ws.Range("R158").Value = "Hortaliza"

# Populate the new row 159 (Segunda) with fresh weekly data
$ws.Range("A159").Value = 8
$ws.Range("B159").Value = "Terminal La Palmera de La Serena"
$ws.Range("C159").Value = "Coquimbo"
$ws.Range("D159").Value = 44554
$ws.Range("E159").Value = 4
$ws.Range("F159").Value = 100112017
$ws.Range("G159").Value = "Apio"
$ws.Range("H159").Value = "Americana (o)"
$ws.Range("I159").Value = "Segunda"
$ws.Range("J159").Value = 1500
$ws.Range("K159").Value = 6000
$ws.Range("L159").Value = 7000
$ws.Range("M159").Value = 6500
$ws.Range("N159").Value = "`$/docena de matas"
$ws.Range("O159").Value = "Provincia del Elquí"
$ws.Range("P159").Value = 1083
$ws.Range("Q159").Value = 6
$ws.Range("R159").Value = "Hortaliza"
